# Botium-utterances scripting edit:
#  - negations / utterance resolution / parametrization
#  - "Hi" / "nice day" utterances replaced by parametrized GREETING_NAME flow
#  - new WHERE_IS_RESTAURANT utterance block added
#  - active sheet/tab & selections swapped back to "Dialogs"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Dialogs"
$ws2 = $wb.Worksheets.Item(2)   # "Utterances"

# ---------------------------------------------------------------------------
# Sheet "Dialogs"
# ---------------------------------------------------------------------------

# A2 was "Hi" (REFCODE-free greeting bot name cell) -> now "GREETING_NAME Bot",
# and picks up the wrap-text-only style (no special text numberformat).
$ws1.Range("A2").Value = "GREETING_NAME Bot"
$ws1.Range("A2").Style = "Normal"
$ws1.Range("A2").WrapText = $true

# A4 was "where is the next restaurant" -> becomes refcode "WHERE_IS_RESTAURANT"
# and loses its dedicated style entirely (back to plain default formatting).
$ws1.Range("A4").Value = "WHERE_IS_RESTAURANT"
$ws1.Range("A4").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "Utterances"
# ---------------------------------------------------------------------------

# The old "OK" refcode row (row 5) moves up to row 4, and "nice day" becomes "ok".
$ws2.Range("A4").Value = "OK"
$ws2.Range("B4").Value = "ok"
$ws2.Range("A5").ClearContents()
$ws2.Range("B5").Value = "fine"
$ws2.Range("B6").Value = "super"

# New parametrized GREETING_NAME utterance block (rows 7-8).
$ws2.Range("A7").Value = "GREETING_NAME"
$ws2.Range("B7").Value = "hi, %s"
$ws2.Range("B8").Value = "hello, %s"

# New WHERE_IS_RESTAURANT utterance block (rows 9-10), text-formatted like
# the other plain utterance cells that use the dedicated "@" number format.
$ws2.Range("A9").Value = "WHERE_IS_RESTAURANT"
$ws2.Range("B9").Value = "where is the next restaurant"
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B10").Value = "where is a restaurant"
$ws2.Range("B10").NumberFormat = "@"

# Column widths on "Utterances" (new explicit <cols>). The ColumnWidth
# setter only has pixel (1/6-character) granularity, so these are the
# closest achievable values to the authored widths of 24 and 22.5546875
# characters (24 lands exactly; 22.5546875 rounds to the nearest 1/6th).
$ws2.Columns.Item(1).ColumnWidth = 23.16
$ws2.Columns.Item(2).ColumnWidth = 21.7

# ---------------------------------------------------------------------------
# View / selection state: active tab moves back to "Dialogs".
# ---------------------------------------------------------------------------
$ws2.Range("A7").Select()
$ws1.Range("A2").Select()
$ws1.Activate()
